# DiccionarioDatosUD.xlsx - add RENDIMIENTO_* data-dictionary entries (10 new
# "semestre" fields) as rows 349-358 of Sheet1, following the existing
# CODIGO / DESCRIPCION table layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in the new CODIGO / DESCRIPCION rows -------------------------
# Values are entered in the same order the original author typed them in
# (column A first for a couple of rows, then the first description, then
# the rest of column A, then the remaining descriptions) so the shared
# string table comes out in the same sequence.
$ws.Range("A349").Value = "RENDIMIENTO_UNO"
$ws.Range("A350").Value = "RENDIMIENTO_DOS"
$ws.Range("B349").Value = "RENDIMIENTO ESTUDIANTE (1 SEMESTRE)"
$ws.Range("A351").Value = "RENDIMIENTO_TRES"
$ws.Range("A352").Value = "RENDIMIENTO_CUATRO"
$ws.Range("A353").Value = "RENDIMIENTO_CINCO"
$ws.Range("A354").Value = "RENDIMIENTO_SEIS"
$ws.Range("A355").Value = "RENDIMIENTO_SIETE"
$ws.Range("A356").Value = "RENDIMIENTO_OCHO"
$ws.Range("A357").Value = "RENDIMIENTO_NUEVE"
$ws.Range("A358").Value = "RENDIMIENTO_DIEZ"
$ws.Range("B350").Value = "RENDIMIENTO ESTUDIANTE (2 SEMESTRE)"
$ws.Range("B351").Value = "RENDIMIENTO ESTUDIANTE (3 SEMESTRE)"
$ws.Range("B352").Value = "RENDIMIENTO ESTUDIANTE (4 SEMESTRE)"
$ws.Range("B353").Value = "RENDIMIENTO ESTUDIANTE (5 SEMESTRE)"
$ws.Range("B354").Value = "RENDIMIENTO ESTUDIANTE (6 SEMESTRE)"
$ws.Range("B355").Value = "RENDIMIENTO ESTUDIANTE (7 SEMESTRE)"
$ws.Range("B356").Value = "RENDIMIENTO ESTUDIANTE (8 SEMESTRE)"
$ws.Range("B357").Value = "RENDIMIENTO ESTUDIANTE (9 SEMESTRE)"
$ws.Range("B358").Value = "RENDIMIENTO ESTUDIANTE (10 SEMESTRE)"

# --- 2. Formatting ----------------------------------------------------------
# Column B: reuse the exact boxed-border formatting already used throughout
# the table (same as every other DESCRIPCION cell, e.g. B348).
$ws.Range("B348").Copy()
$ws.Range("B349:B358").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A: make sure the new rows get a full thin box border (matching the
# rest of the CODIGO column) instead of the old "left/right only" border
# that used to sit on the former last (blank) row.
$rngA = $ws.Range("A349:A358")
$rngA.Borders.LineStyle = -4142
$rngA.Borders.LineStyle = 1

# --- 3. Update the view / selection -----------------------------------------
$excel.ActiveWindow.ScrollRow = 347
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C360").Select()
